$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.225560903549194
$ws.Range("B1").Value = 2.576617956161499
$ws.Range("C1").Value = 9.238280296325684
$ws.Range("D1").Value = 2.052587985992432
$ws.Range("E1").Value = 1.180825114250183
